$wb = $excel.ActiveWorkbook

# --- Master sheet: add two new test-case rows ---
$master = $wb.Worksheets.Item("Master")

# Clone the formatting (style) of the existing data row onto the two new
# rows before filling in values, so the new cells reuse the same cellXf
# (border-only style) instead of Excel minting a brand-new style entry.
$master.Range("A2:E2").Copy()
$master.Range("A3:E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$master.Range("A3").Value = 2
$master.Range("B3").Value = "MyInfo"
$master.Range("C3").Value = "Validate admin Section"
$master.Range("D3").Value = "LoginPage:Manager;"

$master.Range("A4").Value = 3
$master.Range("B4").Value = "MyInfo"
$master.Range("C4").Value = "Verify tab on Side Search bar"
$master.Range("D4").Value = "LoginPage:Manager;"

# Update selection to match the author's final cursor position
$master.Range("C4").Select() | Out-Null

# --- LoginPage sheet: selection change only ---
$loginPage = $wb.Worksheets.Item("LoginPage")
$loginPage.Range("C1").Select() | Out-Null

$master.Activate() | Out-Null
